$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new exam semester row (2023 - Vaar) at the bottom of the table
$ws.Range("A14").Value = "2023 - Vår"
$ws.Range("B14").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/hjemme-23-v.pdf)"
$ws.Range("C14").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/hjemme-23-v-solprop.html)"
$ws.Range("D14").Value = "[Materiale](tidligere-eksamensoppgaver/hjemme-23-v-ekstra.zip)"

# Match the selection state left behind in the saved workbook
$ws.Range("A15").Select()
